$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 17 de Julio de 2020 a las 08:30"

# Row 6 - India
$ws.Range("B6").Value = 1005760
$ws.Range("C6").Value = 123
$ws.Range("D6").Value = 636660
$ws.Range("E6").Value = 343481
$ws.Range("G6").Value = 10
$ws.Range("H6").Value = 25619

# Row 38 - Ucrania
$ws.Range("B38").Value = 57264
$ws.Range("C38").Value = 809
$ws.Range("D38").Value = 29769
$ws.Range("E38").Value = 26039
$ws.Range("G38").Value = 11
$ws.Range("H38").Value = 1456

# Row 67 - Uzbekistan
$ws.Range("B67").Value = 15349
$ws.Range("C67").Value = 283
$ws.Range("E67").Value = 6490
$ws.Range("G67").Value = 1
$ws.Range("H67").Value = 76

# Row 75 - El Salvador
$ws.Range("D75").Value = 6281
$ws.Range("E75").Value = 4367
$ws.Range("G75").Value = 11
$ws.Range("H75").Value = 309

# Row 146 - Georgia
$ws.Range("B146").Value = 1010
$ws.Range("C146").Value = 4
$ws.Range("D146").Value = 885
$ws.Range("E146").Value = 110
